$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 80, shifting existing rows 80-115 down to 81-116
$ws.Rows.Item(80).Insert()

# Populate new row 80 with data
$ws.Range("A80").Value = 11
$ws.Range("B80").Value = "Vega Monumental Concepción"
$ws.Range("C80").Value = "Bíobío"
$ws.Range("D80").Value = 44992
$ws.Range("E80").Value = 8
$ws.Range("F80").Value = "Fruta"
$ws.Range("G80").Value = 100103
$ws.Range("H80").Value = "Frutos de hueso (carozo)"
$ws.Range("I80").Value = 100103002
$ws.Range("J80").Value = "Ciruela"
$ws.Range("K80").Value = "Larry Ann"
$ws.Range("L80").Value = "Primera"
$ws.Range("M80").Value = 220
$ws.Range("N80").Value = 9000
$ws.Range("O80").Value = 9500
$ws.Range("P80").Value = 9227
$ws.Range("Q80").Value = "$/bandeja 18 kilos granel"
$ws.Range("R80").Value = "Región de O'Higgins"
$ws.Range("S80").Value = 513
$ws.Range("T80").Value = 18
